# Natmi following Dr Hou advice:
# Expand the Lgi2 -> Adam23 sending/target cluster matrix from
# {FAPs, sCs} x {ECs, FAPs, sCs} (6 rows) to the full
# {ECs, FAPs, sCs} x {ECs, FAPs, sCs} (9 rows) grid, and refresh every
# numeric column (ligand/receptor expression + specificity + edge weights)
# with the recomputed values for the expanded cluster set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("ECs","Lgi2","Adam23","ECs",[double]"1.0",[double]"0.3333333333333333",[double]"0.058936",[double]"0.176808",[double]"0.003640699631737656",[double]"0.003640699631737656",[double]"3.0",[double]"1.0",[double]"0.1145763333333333",[double]"0.343729",[double]"0.006557053879060051",[double]"0.006557053879060051",[double]"0.006752670781333332",[double]"0.06077403703199999",[double]"2.38722636427779e-05",[double]"2.38722636427779e-05"),
    @("ECs","Lgi2","Adam23","FAPs",[double]"1.0",[double]"0.3333333333333333",[double]"0.058936",[double]"0.176808",[double]"0.003640699631737656",[double]"0.003640699631737656",[double]"3.0",[double]"1.0",[double]"9.390663666666667",[double]"28.171991",[double]"0.5374154140831726",[double]"0.5374154140831726",[double]"0.5534481538586666",[double]"4.981033384728",[double]"0.001956568100142747",[double]"0.001956568100142747"),
    @("ECs","Lgi2","Adam23","sCs",[double]"1.0",[double]"0.3333333333333333",[double]"0.058936",[double]"0.176808",[double]"0.003640699631737656",[double]"0.003640699631737656",[double]"3.0",[double]"1.0",[double]"7.968512",[double]"23.905536",[double]"0.4560275320377672",[double]"0.4560275320377672",[double]"0.4696322232319999",[double]"4.226690009087999",[double]"0.001660259267952131",[double]"0.001660259267952131"),
    @("FAPs","Lgi2","Adam23","ECs",[double]"3.0",[double]"1.0",[double]"14.32797533333333",[double]"42.983926",[double]"0.8850932285803735",[double]"0.8850932285803734",[double]"3.0",[double]"1.0",[double]"0.1145763333333333",[double]"0.343729",[double]"0.006557053879060051",[double]"0.006557053879060051",[double]"1.641646877783778",[double]"14.774821900054",[double]"0.005803603987792723",[double]"0.005803603987792722"),
    @("FAPs","Lgi2","Adam23","FAPs",[double]"3.0",[double]"1.0",[double]"14.32797533333333",[double]"42.983926",[double]"0.8850932285803735",[double]"0.8850932285803734",[double]"3.0",[double]"1.0",[double]"9.390663666666667",[double]"28.171991",[double]"0.5374154140831726",[double]"0.5374154140831726",[double]"134.5491973796296",[double]"1210.942776416666",[double]"0.4756627439397336",[double]"0.4756627439397335"),
    @("FAPs","Lgi2","Adam23","sCs",[double]"3.0",[double]"1.0",[double]"14.32797533333333",[double]"42.983926",[double]"0.8850932285803735",[double]"0.8850932285803734",[double]"3.0",[double]"1.0",[double]"7.968512",[double]"23.905536",[double]"0.4560275320377672",[double]"0.4560275320377672",[double]"114.1726433793707",[double]"1027.553790414336",[double]"0.4036268806528471",[double]"0.4036268806528471"),
    @("sCs","Lgi2","Adam23","ECs",[double]"3.0",[double]"1.0",[double]"1.801186",[double]"5.403558",[double]"0.111266071787889",[double]"0.1112660717878889",[double]"3.0",[double]"1.0",[double]"0.1145763333333333",[double]"0.343729",[double]"0.006557053879060051",[double]"0.006557053879060051",[double]"0.2063732875313333",[double]"1.857359587782",[double]"0.0007295776276245513",[double]"0.0007295776276245512"),
    @("sCs","Lgi2","Adam23","FAPs",[double]"3.0",[double]"1.0",[double]"1.801186",[double]"5.403558",[double]"0.111266071787889",[double]"0.1112660717878889",[double]"3.0",[double]"1.0",[double]"9.390663666666667",[double]"28.171991",[double]"0.5374154140831726",[double]"0.5374154140831726",[double]"16.91433192710867",[double]"152.228987343978",[double]"0.05979610204329636",[double]"0.05979610204329635"),
    @("sCs","Lgi2","Adam23","sCs",[double]"3.0",[double]"1.0",[double]"1.801186",[double]"5.403558",[double]"0.111266071787889",[double]"0.1112660717878889",[double]"3.0",[double]"1.0",[double]"7.968512",[double]"23.905536",[double]"0.4560275320377672",[double]"0.4560275320377672",[double]"14.352772255232",[double]"129.174950297088",[double]"0.05074039211696804",[double]"0.05074039211696803")
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $rowsData[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
